# Apply analysis dashboard updates (battery-wise) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.02760221064814815
$ws.Range("B2").Value = 20.89532472222222
$ws.Range("B3").Value = 1057.992809768889
$ws.Range("B5").Value = 3.429
$ws.Range("B6").Value = 63
$ws.Range("B7").Value = 8
$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("B8").Value = 24.51783489577219
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("B9").Value = 43.15196730325186
$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("B11").Value = "Custom mode`n57.25%`nEco mode`n42.75%"
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("B13").Value = -1602.345021105595
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.005770237904991706
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.338
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.034
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 7
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 53
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.107452210277778
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.0000001290375897508596
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 1.060945971209572
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 4.673770798280053
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 6.566647971583474
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 12.73135165451486
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 52.18265096279679
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 9.894372779958871
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 7.393905402879043
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 4.40269209197981
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 0.9768180968405309
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
